$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Set Price column (D2:D51) to text format so numeric-looking strings (e.g. "0.9979")
# are preserved exactly as typed, then restore the original (default) style afterwards.
$ws.Range("D2:D51").NumberFormat = "@"

# Update Price (D) column values
$ws.Range("D2").Value = "28.096.87"
$ws.Range("D3").Value = "1.816.82"
$ws.Range("D4").Value = "0.9979"
$ws.Range("D5").Value = "338.26"
$ws.Range("D6").Value = "0.9957"
$ws.Range("D7").Value = "0.3930"
$ws.Range("D8").Value = "0.3490"
$ws.Range("D9").Value = "48.48"
$ws.Range("D10").Value = "1.201"
$ws.Range("D11").Value = "0.07592"
$ws.Range("D12").Value = "0.9966"
$ws.Range("D14").Value = "6.536"
$ws.Range("D15").Value = "1.817.98"
$ws.Range("D16").Value = "7.196"
$ws.Range("D17").Value = "0.00001107"
$ws.Range("D18").Value = "0.06717"
$ws.Range("D19").Value = "85.30"
$ws.Range("D20").Value = "0.9968"
$ws.Range("D21").Value = "17.86"
$ws.Range("D22").Value = "6.572"
$ws.Range("D23").Value = "28.060.73"
$ws.Range("D24").Value = "12.83"
$ws.Range("D25").Value = "2.405"
$ws.Range("D26").Value = "1.516"
$ws.Range("D27").Value = "2.568"
$ws.Range("D28").Value = "21.40"
$ws.Range("D29").Value = "154.73"
$ws.Range("D30").Value = "2.019.00"
$ws.Range("D31").Value = "135.53"
$ws.Range("D32").Value = "4.038"
$ws.Range("D33").Value = "6.138"
$ws.Range("D34").Value = "0.08851"
$ws.Range("D35").Value = "13.31"
$ws.Range("D36").Value = "5.530"
$ws.Range("D37").Value = "0.6953"
$ws.Range("D38").Value = "0.02430"
$ws.Range("D39").Value = "0.06556"
$ws.Range("D40").Value = "1.617"
$ws.Range("D41").Value = "0.2231"
$ws.Range("D42").Value = "1.270"
$ws.Range("D43").Value = "8.528"
$ws.Range("D44").Value = "14.73"
$ws.Range("D45").Value = "0.6526"
$ws.Range("D46").Value = "0.9958"
$ws.Range("D47").Value = "3.876"
$ws.Range("D48").Value = "2.171"
$ws.Range("D49").Value = "132.72"
$ws.Range("D50").Value = "0.07216"
$ws.Range("D51").Value = "80.40"

# Restore default (Normal) styling on the Price column so only the text changed
$ws.Range("D2:D51").Style = "Normal"

# Update Volume(1h) (E) column values
$ws.Range("E2").Value = "  +2.48%  "
$ws.Range("E3").Value = "  +1.14%  "
$ws.Range("E4").Value = "  -0.47%  "
$ws.Range("E5").Value = "  +0.07%  "
$ws.Range("E6").Value = "  -0.55%  "
$ws.Range("E8").Value = "  +0.91%  "
$ws.Range("E9").Value = "  +0.41%  "
$ws.Range("E10").Value = "  -0.43%  "
$ws.Range("E11").Value = "  +0.63%  "
$ws.Range("E12").Value = "  -0.37%  "
$ws.Range("E13").Value = "  +0.24%  "
$ws.Range("E14").Value = "  +0.45%  "
$ws.Range("E15").Value = "  +1.28%  "
$ws.Range("E16").Value = "  +1.29%  "
$ws.Range("E17").Value = "  +0.36%  "
$ws.Range("E18").Value = "  +0.66%  "
$ws.Range("E19").Value = "  +0.12%  "
$ws.Range("E20").Value = "  -0.37%  "
$ws.Range("E21").Value = "  +2.49%  "
$ws.Range("E22").Value = "  +0.55%  "
$ws.Range("E23").Value = "  +2.48%  "
$ws.Range("E24").Value = "  +1.96%  "
$ws.Range("E25").Value = "  -1.27%  "
$ws.Range("E26").Value = "  +0.62%  "
$ws.Range("E27").Value = "  -1.07%  "
$ws.Range("E28").Value = "  -0.34%  "
$ws.Range("E29").Value = "  +1.60%  "
$ws.Range("E30").Value = "  +1.09%  "
$ws.Range("E31").Value = "  +1.02%  "
$ws.Range("E32").Value = "  -0.79%  "
$ws.Range("E33").Value = "  -0.29%  "
$ws.Range("E34").Value = "  +1.37%  "
$ws.Range("E35").Value = "  -0.24%  "
$ws.Range("E36").Value = "  +0.91%  "
$ws.Range("E37").Value = "  +0.26%  "
$ws.Range("E38").Value = "  +3.70%  "
$ws.Range("E39").Value = "  +2.61%  "
$ws.Range("E40").Value = "  -4.38%  "
$ws.Range("E41").Value = "  +0.91%  "
$ws.Range("E42").Value = "  -0.56%  "
$ws.Range("E43").Value = "  -4.42%  "
$ws.Range("E44").Value = "  +2.23%  "
$ws.Range("E45").Value = "  +0.47%  "
$ws.Range("E46").Value = "  -0.49%  "
$ws.Range("E47").Value = "  +0.01%  "
$ws.Range("E48").Value = "  +1.76%  "
$ws.Range("E49").Value = "  +1.59%  "
$ws.Range("E50").Value = "  +0.44%  "
$ws.Range("E51").Value = "  +1.17%  "
